# Deploy 0.614 27.03.2025 - Import barcodes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column D ("of_price") to host the new
#    "of_barcode" column. This shifts of_price/of_img/spare column right
#    by one (D->E, E->F, F->G) while Excel copies the left neighbour's
#    style into the freshly inserted column.
$ws.Columns.Item(4).Insert()

# 2. Header for the new column
$ws.Range("D1").Value = "of_barcode"

# 3. Barcode values (stored as text, matching the "t=s" shared-string
#    cells used for every other text column in the sheet)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2003357847234"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2005847345742,23445542345"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "3453453466544"

# D5 stays empty - "Rubik's 360/Sphere" has no barcode

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3453455667777"

# 4. Column widths: new column D should look like column C.
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# 5. Hyperlinks on the (now shifted) image column F need to be re-pointed;
#    a plain column insert does not relocate them automatically.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://nwms.cloud/upload/catalog/Uzel_2.jpg", "", "", "https://nwms.cloud/upload/catalog/Uzel_2.jpg")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://nwms.cloud/upload/catalog/Izumrud.jpg", "", "", "https://nwms.cloud/upload/catalog/Izumrud.jpg")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://nwms.cloud/upload/catalog/Sudoku-SHar.jpg", "", "", "https://nwms.cloud/upload/catalog/Sudoku-SHar.jpg")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://nwms.cloud/upload/catalog/SHarik_Rubika.jpg", "", "", "https://nwms.cloud/upload/catalog/SHarik_Rubika.jpg")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://nwms.cloud/upload/catalog/YAbloko.jpg", "", "", "https://nwms.cloud/upload/catalog/YAbloko.jpg")

# Adding a hyperlink re-styles the cell with the built-in "Hyperlink" named
# style (different font/xf). Restore the original per-row look (border /
# fill / number format) by pasting the formatting back from column C of
# the same row, which still carries the untouched look.
$ws.Range("C2").Copy() | Out-Null
$ws.Range("F2").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("F3").PasteSpecial(-4122) | Out-Null
$ws.Range("C4").Copy() | Out-Null
$ws.Range("F4").PasteSpecial(-4122) | Out-Null
$ws.Range("C5").Copy() | Out-Null
$ws.Range("F5").PasteSpecial(-4122) | Out-Null
$ws.Range("C6").Copy() | Out-Null
$ws.Range("F6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 6. Add a small decorative footer box (rows 7-10) below the table, drawn
#    as a thin bordered/filled rectangle spanning A7:G10.
$box = $ws.Range("A7:G10")
$box.Interior.ColorIndex = 13
$box.RowHeight = 14.7

$top = $ws.Range("A7:G7")
$top.Borders.Item(8).LineStyle = 1
$top.Borders.Item(8).Color = 0xa5a5a5

$outer = $ws.Range("A7:G10")
$outer.Borders.Item(7).LineStyle = 1
$outer.Borders.Item(7).Color = 0xaaaaaa
$outer.Borders.Item(10).LineStyle = 1
$outer.Borders.Item(10).Color = 0xaaaaaa
$outer.Borders.Item(9).LineStyle = 1
$outer.Borders.Item(9).Color = 0xaaaaaa

Write-Host "done"
